# Updated main GSC export data:
# The oldest day's row ("2025-10-09") is dropped from the rolling-window
# Coverage report on the "Chart" sheet, so every subsequent row shifts up
# by one and the trailing date ("2026-01-03") becomes the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Delete row 2 (the 2025-10-09 entry); everything below shifts up.
$ws.Rows.Item(2).Delete()
